# Auto-generated script applying market-data snapshot updates to each Leve sheet.
# Values correspond to refreshed currentAveragePrice / Leve price / profit columns (H:N).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 140.27272  # H9: 152.5 -> 140.27272
$ws.Cells.Item(9, 9).Value = 130  # I9: 144 -> 130
$ws.Cells.Item(9, 11).Value = 130  # K9: 144 -> 130
$ws.Cells.Item(9, 13).Value = 39  # M9: 25 -> 39
$ws.Cells.Item(12, 10).Value = 956.6667  # J12: 972.5 -> 956.6667
$ws.Cells.Item(12, 12).Value = 956.6667  # L12: 972.5 -> 956.6667
$ws.Cells.Item(12, 14).Value = -1296.6667  # N12: -1312.5 -> -1296.6667
$ws.Cells.Item(32, 8).Value = 2816.1667  # H32: 2784.8572 -> 2816.1667
$ws.Cells.Item(32, 9).Value = 3498.5  # I32: 2748.5 -> 3498.5
$ws.Cells.Item(32, 10).Value = 2475  # J32: 2833.3333 -> 2475
$ws.Cells.Item(32, 11).Value = 3498.5  # K32: 2748.5 -> 3498.5
$ws.Cells.Item(32, 12).Value = 2475  # L32: 2833.3333 -> 2475
$ws.Cells.Item(32, 13).Value = -3172.5  # M32: -2422.5 -> -3172.5
$ws.Cells.Item(32, 14).Value = -3127  # N32: -3485.3333 -> -3127
$ws.Cells.Item(33, 8).Value = 142.8125  # H33: 145.66667 -> 142.8125
$ws.Cells.Item(33, 10).Value = 100  # J33: 0 -> 100
$ws.Cells.Item(33, 12).Value = 100  # L33: 0 -> 100
$ws.Cells.Item(33, 14).Value = -558  # N33: None -> -558
$ws.Cells.Item(39, 8).Value = 146.33333  # H39: 121 -> 146.33333
$ws.Cells.Item(39, 9).Value = 146.33333  # I39: 132.2 -> 146.33333
$ws.Cells.Item(39, 10).Value = 0  # J39: 9 -> 0
$ws.Cells.Item(39, 11).Value = 438.99999  # K39: 396.6 -> 438.99999
$ws.Cells.Item(39, 12).Value = 0  # L39: 27 -> 0
$ws.Cells.Item(39, 13).Value = -142.99999  # M39: -100.6 -> -142.99999
$ws.Cells.Item(39, 14).ClearContents()  # N39: was -619
$ws.Cells.Item(43, 8).Value = 4131.1113  # H43: 4143.3335 -> 4131.1113
$ws.Cells.Item(43, 9).Value = 1946.5  # I43: 1998 -> 1946.5
$ws.Cells.Item(43, 10).Value = 4755.2856  # J43: 4756.2856 -> 4755.2856
$ws.Cells.Item(43, 11).Value = 1946.5  # K43: 1998 -> 1946.5
$ws.Cells.Item(43, 12).Value = 4755.2856  # L43: 4756.2856 -> 4755.2856
$ws.Cells.Item(43, 13).Value = -1877.5  # M43: -1929 -> -1877.5
$ws.Cells.Item(43, 14).Value = -4893.2856  # N43: -4894.2856 -> -4893.2856
$ws.Cells.Item(51, 8).Value = 7966.6665  # H51: 8000 -> 7966.6665
$ws.Cells.Item(51, 10).Value = 7966.6665  # J51: 8000 -> 7966.6665
$ws.Cells.Item(51, 12).Value = 7966.6665  # L51: 8000 -> 7966.6665
$ws.Cells.Item(51, 14).Value = -8934.666499999999  # N51: -8968 -> -8934.666499999999
$ws.Cells.Item(58, 8).Value = 35  # H58: 1923 -> 35
$ws.Cells.Item(58, 9).Value = 35  # I58: 34.5 -> 35
$ws.Cells.Item(58, 10).Value = 0  # J58: 5700 -> 0
$ws.Cells.Item(58, 11).Value = 105  # K58: 103.5 -> 105
$ws.Cells.Item(58, 12).Value = 0  # L58: 17100 -> 0
$ws.Cells.Item(58, 13).Value = 45  # M58: 46.5 -> 45
$ws.Cells.Item(58, 14).ClearContents()  # N58: was -17400
$ws.Cells.Item(64, 8).Value = 4718.1816  # H64: 4445.4546 -> 4718.1816
$ws.Cells.Item(64, 10).Value = 3180  # J64: 3150 -> 3180
$ws.Cells.Item(64, 12).Value = 3180  # L64: 3150 -> 3180
$ws.Cells.Item(64, 14).Value = -3676  # N64: -3646 -> -3676
$ws.Cells.Item(67, 8).Value = 4718.1816  # H67: 4445.4546 -> 4718.1816
$ws.Cells.Item(67, 10).Value = 3180  # J67: 3150 -> 3180
$ws.Cells.Item(67, 12).Value = 3180  # L67: 3150 -> 3180
$ws.Cells.Item(67, 14).Value = -4896  # N67: -4866 -> -4896
$ws.Cells.Item(98, 8).Value = 625.38464  # H98: 628.0769 -> 625.38464
$ws.Cells.Item(98, 9).Value = 625.38464  # I98: 628.0769 -> 625.38464
$ws.Cells.Item(98, 11).Value = 625.38464  # K98: 628.0769 -> 625.38464
$ws.Cells.Item(98, 13).Value = 872.61536  # M98: 869.9231 -> 872.61536
$ws.Cells.Item(106, 8).Value = 34888.223  # H106: 34999.332 -> 34888.223
$ws.Cells.Item(106, 9).Value = 38713.43  # I106: 36249.25 -> 38713.43
$ws.Cells.Item(106, 10).Value = 21500  # J106: 25000 -> 21500
$ws.Cells.Item(106, 11).Value = 38713.43  # K106: 36249.25 -> 38713.43
$ws.Cells.Item(106, 12).Value = 21500  # L106: 25000 -> 21500
$ws.Cells.Item(106, 13).Value = -38082.43  # M106: -35618.25 -> -38082.43
$ws.Cells.Item(106, 14).Value = -22762  # N106: -26262 -> -22762
$ws.Cells.Item(107, 8).Value = 295.1  # H107: 311.22223 -> 295.1
$ws.Cells.Item(107, 9).Value = 295.1  # I107: 311.22223 -> 295.1
$ws.Cells.Item(107, 11).Value = 295.1  # K107: 311.22223 -> 295.1
$ws.Cells.Item(107, 13).Value = 1624.9  # M107: 1608.77777 -> 1624.9
$ws.Cells.Item(122, 8).Value = 625.38464  # H122: 628.0769 -> 625.38464
$ws.Cells.Item(122, 9).Value = 625.38464  # I122: 628.0769 -> 625.38464
$ws.Cells.Item(122, 11).Value = 1876.15392  # K122: 1884.2307 -> 1876.15392
$ws.Cells.Item(122, 13).Value = 573.84608  # M122: 565.7692999999999 -> 573.84608
$ws.Cells.Item(135, 8).Value = 2181.111  # H135: 1986.9 -> 2181.111
$ws.Cells.Item(135, 9).Value = 1733  # I135: 1359.5 -> 1733
$ws.Cells.Item(135, 11).Value = 15597  # K135: 12235.5 -> 15597
$ws.Cells.Item(135, 13).Value = -13062  # M135: -9700.5 -> -13062

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 7838.61  # H32: 7838.6343 -> 7838.61
$ws.Cells.Item(32, 9).Value = 5983.324  # I32: 5983.3516 -> 5983.324
$ws.Cells.Item(32, 11).Value = 5983.324  # K32: 5983.3516 -> 5983.324
$ws.Cells.Item(32, 13).Value = -5696.324  # M32: -5696.3516 -> -5696.324
$ws.Cells.Item(61, 8).Value = 0  # H61: 3851.8572 -> 0
$ws.Cells.Item(61, 9).Value = 0  # I61: 3851.8572 -> 0
$ws.Cells.Item(61, 11).Value = 0  # K61: 3851.8572 -> 0
$ws.Cells.Item(61, 13).ClearContents()  # M61: was -3639.8572
$ws.Cells.Item(63, 8).Value = 6045.5386  # H63: 5741.9287 -> 6045.5386
$ws.Cells.Item(63, 10).Value = 7661.625  # J63: 7009.778 -> 7661.625
$ws.Cells.Item(63, 12).Value = 7661.625  # L63: 7009.778 -> 7661.625
$ws.Cells.Item(63, 14).Value = -9033.625  # N63: -8381.778 -> -9033.625
$ws.Cells.Item(66, 8).Value = 6045.5386  # H66: 5741.9287 -> 6045.5386
$ws.Cells.Item(66, 10).Value = 7661.625  # J66: 7009.778 -> 7661.625
$ws.Cells.Item(66, 12).Value = 38308.125  # L66: 35048.89 -> 38308.125
$ws.Cells.Item(66, 14).Value = -45172.125  # N66: -41912.89 -> -45172.125
$ws.Cells.Item(97, 8).Value = 2501.3157  # H97: 2618.6667 -> 2501.3157
$ws.Cells.Item(97, 9).Value = 1447.9333  # I97: 1523.5714 -> 1447.9333
$ws.Cells.Item(97, 11).Value = 1447.9333  # K97: 1523.5714 -> 1447.9333
$ws.Cells.Item(97, 13).Value = -951.9332999999999  # M97: -1027.5714 -> -951.9332999999999
$ws.Cells.Item(110, 8).Value = 4891.5713  # H110: 8074.25 -> 4891.5713
$ws.Cells.Item(110, 9).Value = 5599.3335  # I110: 8074.25 -> 5599.3335
$ws.Cells.Item(110, 10).Value = 645  # J110: 0 -> 645
$ws.Cells.Item(110, 11).Value = 5599.3335  # K110: 8074.25 -> 5599.3335
$ws.Cells.Item(110, 12).Value = 645  # L110: 0 -> 645
$ws.Cells.Item(110, 13).Value = -3554.3335  # M110: -6029.25 -> -3554.3335
$ws.Cells.Item(110, 14).Value = -4735  # N110: None -> -4735
$ws.Cells.Item(132, 8).Value = 2191.2  # H132: 2253.6667 -> 2191.2
$ws.Cells.Item(132, 9).Value = 1853.6818  # I132: 1909 -> 1853.6818
$ws.Cells.Item(132, 11).Value = 5561.0454  # K132: 5727 -> 5561.0454
$ws.Cells.Item(132, 13).Value = -3031.0454  # M132: -3197 -> -3031.0454
$ws.Cells.Item(136, 8).Value = 0  # H136: 3851.8572 -> 0
$ws.Cells.Item(136, 9).Value = 0  # I136: 3851.8572 -> 0
$ws.Cells.Item(136, 11).Value = 0  # K136: 11555.5716 -> 0
$ws.Cells.Item(136, 13).ClearContents()  # M136: was -9005.571599999999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 947.52  # H94: 915.8148 -> 947.52
$ws.Cells.Item(94, 9).Value = 975.56525  # I94: 939.08 -> 975.56525
$ws.Cells.Item(94, 11).Value = 975.56525  # K94: 939.08 -> 975.56525
$ws.Cells.Item(94, 13).Value = -524.56525  # M94: -488.08 -> -524.56525
$ws.Cells.Item(132, 8).Value = 100584.75  # H132: 104556.164 -> 100584.75
$ws.Cells.Item(132, 10).Value = 100584.75  # J132: 104556.164 -> 100584.75
$ws.Cells.Item(132, 12).Value = 100584.75  # L132: 104556.164 -> 100584.75
$ws.Cells.Item(132, 14).Value = -110704.75  # N132: -114676.164 -> -110704.75
$ws.Cells.Item(133, 8).Value = 150000  # H133: 0 -> 150000
$ws.Cells.Item(133, 10).Value = 150000  # J133: 0 -> 150000
$ws.Cells.Item(133, 12).Value = 150000  # L133: 0 -> 150000
$ws.Cells.Item(133, 14).Value = -160120  # N133: None -> -160120

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 745.4286  # H16: 814.5454999999999 -> 745.4286
$ws.Cells.Item(16, 9).Value = 757.8889  # I16: 835 -> 757.8889
$ws.Cells.Item(16, 10).Value = 723  # J16: 778.75 -> 723
$ws.Cells.Item(16, 11).Value = 757.8889  # K16: 835 -> 757.8889
$ws.Cells.Item(16, 12).Value = 723  # L16: 778.75 -> 723
$ws.Cells.Item(16, 13).Value = -470.8889  # M16: -548 -> -470.8889
$ws.Cells.Item(16, 14).Value = -1297  # N16: -1352.75 -> -1297
$ws.Cells.Item(62, 8).Value = 73898.664  # H62: 87558.39999999999 -> 73898.664
$ws.Cells.Item(62, 10).Value = 73898.664  # J62: 87558.39999999999 -> 73898.664
$ws.Cells.Item(62, 12).Value = 73898.664  # L62: 87558.39999999999 -> 73898.664
$ws.Cells.Item(62, 14).Value = -75146.664  # N62: -88806.39999999999 -> -75146.664
$ws.Cells.Item(65, 8).Value = 73898.664  # H65: 87558.39999999999 -> 73898.664
$ws.Cells.Item(65, 10).Value = 73898.664  # J65: 87558.39999999999 -> 73898.664
$ws.Cells.Item(65, 12).Value = 369493.32  # L65: 437792 -> 369493.32
$ws.Cells.Item(65, 14).Value = -375733.32  # N65: -444032 -> -375733.32
$ws.Cells.Item(105, 8).Value = 496.55554  # H105: 446.4 -> 496.55554
$ws.Cells.Item(105, 9).Value = 508.17648  # I105: 446.4 -> 508.17648
$ws.Cells.Item(105, 10).Value = 299  # J105: 0 -> 299
$ws.Cells.Item(105, 11).Value = 508.17648  # K105: 446.4 -> 508.17648
$ws.Cells.Item(105, 12).Value = 299  # L105: 0 -> 299
$ws.Cells.Item(105, 13).Value = 1238.82352  # M105: 1300.6 -> 1238.82352
$ws.Cells.Item(105, 14).Value = -3793  # N105: None -> -3793
$ws.Cells.Item(113, 8).Value = 745.4286  # H113: 814.5454999999999 -> 745.4286
$ws.Cells.Item(113, 9).Value = 757.8889  # I113: 835 -> 757.8889
$ws.Cells.Item(113, 10).Value = 723  # J113: 778.75 -> 723
$ws.Cells.Item(113, 11).Value = 757.8889  # K113: 835 -> 757.8889
$ws.Cells.Item(113, 12).Value = 723  # L113: 778.75 -> 723
$ws.Cells.Item(113, 13).Value = 1412.1111  # M113: 1335 -> 1412.1111
$ws.Cells.Item(113, 14).Value = -5063  # N113: -5118.75 -> -5063

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 2601086.8  # H4: 2415305.8 -> 2601086.8
$ws.Cells.Item(4, 10).Value = 3550  # J4: 2700 -> 3550
$ws.Cells.Item(4, 12).Value = 10650  # L4: 8100 -> 10650
$ws.Cells.Item(4, 14).Value = -10874  # N4: -8324 -> -10874

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 2282.0833  # H132: 2148.8572 -> 2282.0833
$ws.Cells.Item(132, 9).Value = 1321.8889  # I132: 1326.909 -> 1321.8889
$ws.Cells.Item(132, 11).Value = 3965.6667  # K132: 3980.727 -> 3965.6667
$ws.Cells.Item(132, 13).Value = -1435.6667  # M132: -1450.727 -> -1435.6667

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 3982.6667  # H7: 3949 -> 3982.6667
$ws.Cells.Item(7, 9).Value = 3474  # I7: 3949 -> 3474
$ws.Cells.Item(7, 10).Value = 5000  # J7: 0 -> 5000
$ws.Cells.Item(7, 11).Value = 3474  # K7: 3949 -> 3474
$ws.Cells.Item(7, 12).Value = 5000  # L7: 0 -> 5000
$ws.Cells.Item(7, 13).Value = -3362  # M7: -3837 -> -3362
$ws.Cells.Item(7, 14).Value = -5224  # N7: None -> -5224
$ws.Cells.Item(40, 8).Value = 1844.8  # H40: 1858.4 -> 1844.8
$ws.Cells.Item(40, 9).Value = 1844.8  # I40: 1848 -> 1844.8
$ws.Cells.Item(40, 10).Value = 0  # J40: 1900 -> 0
$ws.Cells.Item(40, 11).Value = 1844.8  # K40: 1848 -> 1844.8
$ws.Cells.Item(40, 12).Value = 0  # L40: 1900 -> 0
$ws.Cells.Item(40, 13).Value = -1708.8  # M40: -1712 -> -1708.8
$ws.Cells.Item(40, 14).ClearContents()  # N40: was -2172
$ws.Cells.Item(43, 8).Value = 10333.333  # H43: 10400 -> 10333.333
$ws.Cells.Item(43, 10).Value = 10333.333  # J43: 10400 -> 10333.333
$ws.Cells.Item(43, 12).Value = 10333.333  # L43: 10400 -> 10333.333
$ws.Cells.Item(43, 14).Value = -10719.333  # N43: -10786 -> -10719.333
$ws.Cells.Item(46, 8).Value = 1099.6666  # H46: 1013.2857 -> 1099.6666
$ws.Cells.Item(46, 9).Value = 999.3333  # I46: 873.25 -> 999.3333
$ws.Cells.Item(46, 11).Value = 999.3333  # K46: 873.25 -> 999.3333
$ws.Cells.Item(46, 13).Value = -811.3333  # M46: -685.25 -> -811.3333
$ws.Cells.Item(82, 8).Value = 3409  # H82: 3234.9473 -> 3409
$ws.Cells.Item(82, 9).Value = 3670.5454  # I82: 3364.8333 -> 3670.5454
$ws.Cells.Item(82, 10).Value = 2998  # J82: 3012.2856 -> 2998
$ws.Cells.Item(82, 11).Value = 3670.5454  # K82: 3364.8333 -> 3670.5454
$ws.Cells.Item(82, 12).Value = 2998  # L82: 3012.2856 -> 2998
$ws.Cells.Item(82, 13).Value = -3309.5454  # M82: -3003.8333 -> -3309.5454
$ws.Cells.Item(82, 14).Value = -3720  # N82: -3734.2856 -> -3720
$ws.Cells.Item(85, 8).Value = 3409  # H85: 3234.9473 -> 3409
$ws.Cells.Item(85, 9).Value = 3670.5454  # I85: 3364.8333 -> 3670.5454
$ws.Cells.Item(85, 10).Value = 2998  # J85: 3012.2856 -> 2998
$ws.Cells.Item(85, 11).Value = 3670.5454  # K85: 3364.8333 -> 3670.5454
$ws.Cells.Item(85, 12).Value = 2998  # L85: 3012.2856 -> 2998
$ws.Cells.Item(85, 13).Value = -2422.5454  # M85: -2116.8333 -> -2422.5454
$ws.Cells.Item(85, 14).Value = -5494  # N85: -5508.2856 -> -5494
$ws.Cells.Item(126, 8).Value = 3982.6667  # H126: 3949 -> 3982.6667
$ws.Cells.Item(126, 9).Value = 3474  # I126: 3949 -> 3474
$ws.Cells.Item(126, 10).Value = 5000  # J126: 0 -> 5000
$ws.Cells.Item(126, 11).Value = 10422  # K126: 11847 -> 10422
$ws.Cells.Item(126, 12).Value = 15000  # L126: 0 -> 15000
$ws.Cells.Item(126, 13).Value = -7952  # M126: -9377 -> -7952
$ws.Cells.Item(126, 14).Value = -19940  # N126: None -> -19940
$ws.Cells.Item(132, 8).Value = 5003  # H132: 0 -> 5003
$ws.Cells.Item(132, 9).Value = 5003  # I132: 0 -> 5003
$ws.Cells.Item(132, 11).Value = 15009  # K132: 0 -> 15009
$ws.Cells.Item(132, 13).Value = -12479  # M132: None -> -12479
$ws.Cells.Item(136, 8).Value = 4000  # H136: 0 -> 4000
$ws.Cells.Item(136, 9).Value = 4000  # I136: 0 -> 4000
$ws.Cells.Item(136, 11).Value = 12000  # K136: 0 -> 12000
$ws.Cells.Item(136, 13).Value = -9450  # M136: None -> -9450

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(15, 8).Value = 40007  # H15: 0 -> 40007
$ws.Cells.Item(15, 10).Value = 40007  # J15: 0 -> 40007
$ws.Cells.Item(15, 12).Value = 40007  # L15: 0 -> 40007
$ws.Cells.Item(15, 14).Value = -40583  # N15: None -> -40583
$ws.Cells.Item(81, 8).Value = 2450.8667  # H81: 2501.9333 -> 2450.8667
$ws.Cells.Item(81, 9).Value = 1906.2307  # I81: 2079 -> 1906.2307
$ws.Cells.Item(81, 10).Value = 5991  # J81: 4193.6665 -> 5991
$ws.Cells.Item(81, 11).Value = 3812.4614  # K81: 4158 -> 3812.4614
$ws.Cells.Item(81, 12).Value = 11982  # L81: 8387.333000000001 -> 11982
$ws.Cells.Item(81, 13).Value = -2751.4614  # M81: -3097 -> -2751.4614
$ws.Cells.Item(81, 14).Value = -14104  # N81: -10509.333 -> -14104
$ws.Cells.Item(84, 8).Value = 2450.8667  # H84: 2501.9333 -> 2450.8667
$ws.Cells.Item(84, 9).Value = 1906.2307  # I84: 2079 -> 1906.2307
$ws.Cells.Item(84, 10).Value = 5991  # J84: 4193.6665 -> 5991
$ws.Cells.Item(84, 11).Value = 19062.307  # K84: 20790 -> 19062.307
$ws.Cells.Item(84, 12).Value = 59910  # L84: 41936.665 -> 59910
$ws.Cells.Item(84, 13).Value = -13758.307  # M84: -15486 -> -13758.307
$ws.Cells.Item(84, 14).Value = -70518  # N84: -52544.665 -> -70518
$ws.Cells.Item(107, 8).Value = 593.4  # H107: 614.3570999999999 -> 593.4
$ws.Cells.Item(107, 10).Value = 704.3333  # J107: 754.875 -> 704.3333
$ws.Cells.Item(107, 12).Value = 2112.9999  # L107: 2264.625 -> 2112.9999
$ws.Cells.Item(107, 14).Value = -5952.9999  # N107: -6104.625 -> -5952.9999
$ws.Cells.Item(113, 8).Value = 495.69232  # H113: 527 -> 495.69232
$ws.Cells.Item(113, 9).Value = 583.6667  # I113: 641.625 -> 583.6667
$ws.Cells.Item(113, 11).Value = 1751.0001  # K113: 1924.875 -> 1751.0001
$ws.Cells.Item(113, 13).Value = 418.9999  # M113: 245.125 -> 418.9999
$ws.Cells.Item(126, 8).Value = 2018.125  # H126: 2068.913 -> 2018.125
$ws.Cells.Item(126, 10).Value = 2884.4285  # J126: 3223.5 -> 2884.4285
$ws.Cells.Item(126, 12).Value = 8653.2855  # L126: 9670.5 -> 8653.2855
$ws.Cells.Item(126, 14).Value = -13593.2855  # N126: -14610.5 -> -13593.2855
$ws.Cells.Item(132, 8).Value = 0  # H132: 1817.1666 -> 0
$ws.Cells.Item(132, 9).Value = 0  # I132: 1634.3334 -> 0
$ws.Cells.Item(132, 10).Value = 0  # J132: 2000 -> 0
$ws.Cells.Item(132, 11).Value = 0  # K132: 4903.0002 -> 0
$ws.Cells.Item(132, 12).Value = 0  # L132: 6000 -> 0
$ws.Cells.Item(132, 13).ClearContents()  # M132: was -2373.0002
$ws.Cells.Item(132, 14).ClearContents()  # N132: was -11060
$ws.Cells.Item(136, 8).Value = 7644.1665  # H136: 8193 -> 7644.1665
$ws.Cells.Item(136, 9).Value = 7967.75  # I136: 8990.333000000001 -> 7967.75
$ws.Cells.Item(136, 11).Value = 23903.25  # K136: 26970.999 -> 23903.25
$ws.Cells.Item(136, 13).Value = -21353.25  # M136: -24420.999 -> -21353.25

Write-Host "Applied updates to Leve profit sheets"